# Update latest optimisation result (run 37)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# --- Sheet1 (Schedule): update totals for existing rows 2 and 3 ---
$ws1.Cells.Item(2,5).Value = 930.0173805000001
$ws1.Cells.Item(2,6).Value = 20.50302867063492
$ws1.Cells.Item(3,5).Value = 324.186837
$ws1.Cells.Item(3,6).Value = 21.44092837301587

# --- Sheet1 (Schedule): add new row 4 ---
$ws1.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(4,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(4,1).Value = 46039.3125
$ws1.Cells.Item(4,2).Value = 46039.8125
$ws1.Cells.Item(4,3).Value = 12
$ws1.Cells.Item(4,4).Value = 45.36
$ws1.Cells.Item(4,5).Value = -112.892013
$ws1.Cells.Item(4,6).Value = -2.488800992063492

# --- Sheet2 (Detailed): update existing rows 25-47 ---
$ws2.Cells.Item(25,2).Value = 36.05989
$ws2.Cells.Item(26,2).Value = 36.06
$ws2.Cells.Item(27,2).Value = 40.54
$ws2.Cells.Item(27,3).Value = "historical"
$ws2.Cells.Item(28,2).Value = 36.0601
$ws2.Cells.Item(28,3).Value = "historical"
$ws2.Cells.Item(29,2).Value = -22.255
$ws2.Cells.Item(29,3).Value = "historical"
$ws2.Cells.Item(30,2).Value = -14
$ws2.Cells.Item(30,3).Value = "historical"
$ws2.Cells.Item(31,2).Value = -14
$ws2.Cells.Item(32,2).Value = -15.63074
$ws2.Cells.Item(33,2).Value = -4.20662
$ws2.Cells.Item(34,2).Value = 2.44544
$ws2.Cells.Item(35,2).Value = -4.10876
$ws2.Cells.Item(36,2).Value = -6.73245
$ws2.Cells.Item(37,2).Value = -6.60018
$ws2.Cells.Item(38,2).Value = -3.12186
$ws2.Cells.Item(39,2).Value = -2.73786
$ws2.Cells.Item(40,2).Value = -2.74631
$ws2.Cells.Item(41,2).Value = 12.41698
$ws2.Cells.Item(42,2).Value = 12.51189
$ws2.Cells.Item(43,2).Value = 20.67595
$ws2.Cells.Item(44,2).Value = 6.37544
$ws2.Cells.Item(45,2).Value = 56.98
$ws2.Cells.Item(47,2).Value = 64.93597

# --- Sheet2 (Detailed): add new rows 50-97 ---
$ws2.Range("A50:A97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("D50:D97").NumberFormat = "YYYY-MM-DD"

$ws2.Cells.Item(50,1).Value = 46039
$ws2.Cells.Item(50,2).Value = 56.98
$ws2.Cells.Item(50,3).Value = "forecast"
$ws2.Cells.Item(50,4).Value = 46039
$ws2.Cells.Item(50,5).Value = "OFF"
$ws2.Cells.Item(51,1).Value = 46039.02083333334
$ws2.Cells.Item(51,2).Value = 56.98
$ws2.Cells.Item(51,3).Value = "forecast"
$ws2.Cells.Item(51,4).Value = 46039
$ws2.Cells.Item(51,5).Value = "OFF"
$ws2.Cells.Item(52,1).Value = 46039.04166666666
$ws2.Cells.Item(52,2).Value = 36.06
$ws2.Cells.Item(52,3).Value = "forecast"
$ws2.Cells.Item(52,4).Value = 46039
$ws2.Cells.Item(52,5).Value = "OFF"
$ws2.Cells.Item(53,1).Value = 46039.0625
$ws2.Cells.Item(53,2).Value = 36.06
$ws2.Cells.Item(53,3).Value = "forecast"
$ws2.Cells.Item(53,4).Value = 46039
$ws2.Cells.Item(53,5).Value = "OFF"
$ws2.Cells.Item(54,1).Value = 46039.08333333334
$ws2.Cells.Item(54,2).Value = 35.87992
$ws2.Cells.Item(54,3).Value = "forecast"
$ws2.Cells.Item(54,4).Value = 46039
$ws2.Cells.Item(54,5).Value = "OFF"
$ws2.Cells.Item(55,1).Value = 46039.10416666666
$ws2.Cells.Item(55,2).Value = 36.06
$ws2.Cells.Item(55,3).Value = "forecast"
$ws2.Cells.Item(55,4).Value = 46039
$ws2.Cells.Item(55,5).Value = "OFF"
$ws2.Cells.Item(56,1).Value = 46039.125
$ws2.Cells.Item(56,2).Value = 36.06
$ws2.Cells.Item(56,3).Value = "forecast"
$ws2.Cells.Item(56,4).Value = 46039
$ws2.Cells.Item(56,5).Value = "OFF"
$ws2.Cells.Item(57,1).Value = 46039.14583333334
$ws2.Cells.Item(57,2).Value = 36.06
$ws2.Cells.Item(57,3).Value = "forecast"
$ws2.Cells.Item(57,4).Value = 46039
$ws2.Cells.Item(57,5).Value = "OFF"
$ws2.Cells.Item(58,1).Value = 46039.16666666666
$ws2.Cells.Item(58,2).Value = 36.06
$ws2.Cells.Item(58,3).Value = "forecast"
$ws2.Cells.Item(58,4).Value = 46039
$ws2.Cells.Item(58,5).Value = "OFF"
$ws2.Cells.Item(59,1).Value = 46039.1875
$ws2.Cells.Item(59,2).Value = 56.98
$ws2.Cells.Item(59,3).Value = "forecast"
$ws2.Cells.Item(59,4).Value = 46039
$ws2.Cells.Item(59,5).Value = "OFF"
$ws2.Cells.Item(60,1).Value = 46039.20833333334
$ws2.Cells.Item(60,2).Value = 56.98
$ws2.Cells.Item(60,3).Value = "forecast"
$ws2.Cells.Item(60,4).Value = 46039
$ws2.Cells.Item(60,5).Value = "OFF"
$ws2.Cells.Item(61,1).Value = 46039.22916666666
$ws2.Cells.Item(61,2).Value = 56.98
$ws2.Cells.Item(61,3).Value = "forecast"
$ws2.Cells.Item(61,4).Value = 46039
$ws2.Cells.Item(61,5).Value = "OFF"
$ws2.Cells.Item(62,1).Value = 46039.25
$ws2.Cells.Item(62,2).Value = 52.585
$ws2.Cells.Item(62,3).Value = "forecast"
$ws2.Cells.Item(62,4).Value = 46039
$ws2.Cells.Item(62,5).Value = "OFF"
$ws2.Cells.Item(63,1).Value = 46039.27083333334
$ws2.Cells.Item(63,2).Value = 36.06
$ws2.Cells.Item(63,3).Value = "forecast"
$ws2.Cells.Item(63,4).Value = 46039
$ws2.Cells.Item(63,5).Value = "OFF"
$ws2.Cells.Item(64,1).Value = 46039.29166666666
$ws2.Cells.Item(64,2).Value = 36.06
$ws2.Cells.Item(64,3).Value = "forecast"
$ws2.Cells.Item(64,4).Value = 46039
$ws2.Cells.Item(64,5).Value = "OFF"
$ws2.Cells.Item(65,1).Value = 46039.3125
$ws2.Cells.Item(65,2).Value = 0.51
$ws2.Cells.Item(65,3).Value = "forecast"
$ws2.Cells.Item(65,4).Value = 46039
$ws2.Cells.Item(65,5).Value = "ON"
$ws2.Cells.Item(66,1).Value = 46039.33333333334
$ws2.Cells.Item(66,2).Value = -5.51
$ws2.Cells.Item(66,3).Value = "forecast"
$ws2.Cells.Item(66,4).Value = 46039
$ws2.Cells.Item(66,5).Value = "ON"
$ws2.Cells.Item(67,1).Value = 46039.35416666666
$ws2.Cells.Item(67,2).Value = -5.78275
$ws2.Cells.Item(67,3).Value = "forecast"
$ws2.Cells.Item(67,4).Value = 46039
$ws2.Cells.Item(67,5).Value = "ON"
$ws2.Cells.Item(68,1).Value = 46039.375
$ws2.Cells.Item(68,2).Value = -4.83962
$ws2.Cells.Item(68,3).Value = "forecast"
$ws2.Cells.Item(68,4).Value = 46039
$ws2.Cells.Item(68,5).Value = "ON"
$ws2.Cells.Item(69,1).Value = 46039.39583333334
$ws2.Cells.Item(69,2).Value = -0.9809099999999999
$ws2.Cells.Item(69,3).Value = "forecast"
$ws2.Cells.Item(69,4).Value = 46039
$ws2.Cells.Item(69,5).Value = "ON"
$ws2.Cells.Item(70,1).Value = 46039.41666666666
$ws2.Cells.Item(70,2).Value = -0.88414
$ws2.Cells.Item(70,3).Value = "forecast"
$ws2.Cells.Item(70,4).Value = 46039
$ws2.Cells.Item(70,5).Value = "ON"
$ws2.Cells.Item(71,1).Value = 46039.4375
$ws2.Cells.Item(71,2).Value = 22.07
$ws2.Cells.Item(71,3).Value = "forecast"
$ws2.Cells.Item(71,4).Value = 46039
$ws2.Cells.Item(71,5).Value = "ON"
$ws2.Cells.Item(72,1).Value = 46039.45833333334
$ws2.Cells.Item(72,2).Value = 0
$ws2.Cells.Item(72,3).Value = "forecast"
$ws2.Cells.Item(72,4).Value = 46039
$ws2.Cells.Item(72,5).Value = "ON"
$ws2.Cells.Item(73,1).Value = 46039.47916666666
$ws2.Cells.Item(73,2).Value = 7.88135
$ws2.Cells.Item(73,3).Value = "forecast"
$ws2.Cells.Item(73,4).Value = 46039
$ws2.Cells.Item(73,5).Value = "ON"
$ws2.Cells.Item(74,1).Value = 46039.5
$ws2.Cells.Item(74,2).Value = 7.21991
$ws2.Cells.Item(74,3).Value = "forecast"
$ws2.Cells.Item(74,4).Value = 46039
$ws2.Cells.Item(74,5).Value = "ON"
$ws2.Cells.Item(75,1).Value = 46039.52083333334
$ws2.Cells.Item(75,2).Value = -0.83992
$ws2.Cells.Item(75,3).Value = "forecast"
$ws2.Cells.Item(75,4).Value = 46039
$ws2.Cells.Item(75,5).Value = "ON"
$ws2.Cells.Item(76,1).Value = 46039.54166666666
$ws2.Cells.Item(76,2).Value = -4.81333
$ws2.Cells.Item(76,3).Value = "forecast"
$ws2.Cells.Item(76,4).Value = 46039
$ws2.Cells.Item(76,5).Value = "ON"
$ws2.Cells.Item(77,1).Value = 46039.5625
$ws2.Cells.Item(77,2).Value = -9.99
$ws2.Cells.Item(77,3).Value = "forecast"
$ws2.Cells.Item(77,4).Value = 46039
$ws2.Cells.Item(77,5).Value = "ON"
$ws2.Cells.Item(78,1).Value = 46039.58333333334
$ws2.Cells.Item(78,2).Value = -18.98278
$ws2.Cells.Item(78,3).Value = "forecast"
$ws2.Cells.Item(78,4).Value = 46039
$ws2.Cells.Item(78,5).Value = "ON"
$ws2.Cells.Item(79,1).Value = 46039.60416666666
$ws2.Cells.Item(79,2).Value = -12.26241
$ws2.Cells.Item(79,3).Value = "forecast"
$ws2.Cells.Item(79,4).Value = 46039
$ws2.Cells.Item(79,5).Value = "ON"
$ws2.Cells.Item(80,1).Value = 46039.625
$ws2.Cells.Item(80,2).Value = -10.70447
$ws2.Cells.Item(80,3).Value = "forecast"
$ws2.Cells.Item(80,4).Value = 46039
$ws2.Cells.Item(80,5).Value = "ON"
$ws2.Cells.Item(81,1).Value = 46039.64583333334
$ws2.Cells.Item(81,2).Value = -12.11173
$ws2.Cells.Item(81,3).Value = "forecast"
$ws2.Cells.Item(81,4).Value = 46039
$ws2.Cells.Item(81,5).Value = "ON"
$ws2.Cells.Item(82,1).Value = 46039.66666666666
$ws2.Cells.Item(82,2).Value = -13.26203
$ws2.Cells.Item(82,3).Value = "forecast"
$ws2.Cells.Item(82,4).Value = 46039
$ws2.Cells.Item(82,5).Value = "ON"
$ws2.Cells.Item(83,1).Value = 46039.6875
$ws2.Cells.Item(83,2).Value = -22.57861
$ws2.Cells.Item(83,3).Value = "forecast"
$ws2.Cells.Item(83,4).Value = 46039
$ws2.Cells.Item(83,5).Value = "ON"
$ws2.Cells.Item(84,1).Value = 46039.70833333334
$ws2.Cells.Item(84,2).Value = -13.43756
$ws2.Cells.Item(84,3).Value = "forecast"
$ws2.Cells.Item(84,4).Value = 46039
$ws2.Cells.Item(84,5).Value = "ON"
$ws2.Cells.Item(85,1).Value = 46039.72916666666
$ws2.Cells.Item(85,2).Value = -9.5
$ws2.Cells.Item(85,3).Value = "forecast"
$ws2.Cells.Item(85,4).Value = 46039
$ws2.Cells.Item(85,5).Value = "ON"
$ws2.Cells.Item(86,1).Value = 46039.75
$ws2.Cells.Item(86,2).Value = -7.16571
$ws2.Cells.Item(86,3).Value = "forecast"
$ws2.Cells.Item(86,4).Value = 46039
$ws2.Cells.Item(86,5).Value = "ON"
$ws2.Cells.Item(87,1).Value = 46039.77083333334
$ws2.Cells.Item(87,2).Value = -2.88098
$ws2.Cells.Item(87,3).Value = "forecast"
$ws2.Cells.Item(87,4).Value = 46039
$ws2.Cells.Item(87,5).Value = "ON"
$ws2.Cells.Item(88,1).Value = 46039.79166666666
$ws2.Cells.Item(88,2).Value = 3.05901
$ws2.Cells.Item(88,3).Value = "forecast"
$ws2.Cells.Item(88,4).Value = 46039
$ws2.Cells.Item(88,5).Value = "ON"
$ws2.Cells.Item(89,1).Value = 46039.8125
$ws2.Cells.Item(89,2).Value = 16.37244
$ws2.Cells.Item(89,3).Value = "forecast"
$ws2.Cells.Item(89,4).Value = 46039
$ws2.Cells.Item(89,5).Value = "OFF"
$ws2.Cells.Item(90,1).Value = 46039.83333333334
$ws2.Cells.Item(90,2).Value = 17.98198
$ws2.Cells.Item(90,3).Value = "forecast"
$ws2.Cells.Item(90,4).Value = 46039
$ws2.Cells.Item(90,5).Value = "OFF"
$ws2.Cells.Item(91,1).Value = 46039.85416666666
$ws2.Cells.Item(91,2).Value = 13.59537
$ws2.Cells.Item(91,3).Value = "forecast"
$ws2.Cells.Item(91,4).Value = 46039
$ws2.Cells.Item(91,5).Value = "OFF"
$ws2.Cells.Item(92,1).Value = 46039.875
$ws2.Cells.Item(92,2).Value = 55.12255
$ws2.Cells.Item(92,3).Value = "forecast"
$ws2.Cells.Item(92,4).Value = 46039
$ws2.Cells.Item(92,5).Value = "OFF"
$ws2.Cells.Item(93,1).Value = 46039.89583333334
$ws2.Cells.Item(93,2).Value = 56.98
$ws2.Cells.Item(93,3).Value = "forecast"
$ws2.Cells.Item(93,4).Value = 46039
$ws2.Cells.Item(93,5).Value = "OFF"
$ws2.Cells.Item(94,1).Value = 46039.91666666666
$ws2.Cells.Item(94,2).Value = 47.26597
$ws2.Cells.Item(94,3).Value = "forecast"
$ws2.Cells.Item(94,4).Value = 46039
$ws2.Cells.Item(94,5).Value = "OFF"
$ws2.Cells.Item(95,1).Value = 46039.9375
$ws2.Cells.Item(95,2).Value = 36.06041
$ws2.Cells.Item(95,3).Value = "forecast"
$ws2.Cells.Item(95,4).Value = 46039
$ws2.Cells.Item(95,5).Value = "OFF"
$ws2.Cells.Item(96,1).Value = 46039.95833333334
$ws2.Cells.Item(96,2).Value = 36.06043
$ws2.Cells.Item(96,3).Value = "forecast"
$ws2.Cells.Item(96,4).Value = 46039
$ws2.Cells.Item(96,5).Value = "OFF"
$ws2.Cells.Item(97,1).Value = 46039.97916666666
$ws2.Cells.Item(97,2).Value = 36.0604
$ws2.Cells.Item(97,3).Value = "forecast"
$ws2.Cells.Item(97,4).Value = 46039
$ws2.Cells.Item(97,5).Value = "OFF"
